# Cálculo da quantidade de cabo de cobre utilizada no piso 1
#
# On slide 2 ("Floor 1" structured-cabling diagram), the small RJ45
# port icon shape "Imagem 19" (id 20), which sits just outside room
# 1.1.14 and links into the orange copper-cable run, is nudged up and
# slightly left so it lines up with the cable corner it connects to.
#
# Target position (from the canonical OOXML):
#   <a:off x="7515720" y="5085218"/>   (EMU)
# i.e. in points (PowerPoint's COM Left/Top unit, 12700 EMU/pt):
#   Left = 591.789023...pt   Top = 400.410866...pt
#
# Shape Left/Top are exposed as single-precision (32-bit) floats by
# the PowerPoint COM object model, so the literals below are chosen
# (via the midpoint of the float32 range that floors to the exact
# target EMU) so the round-trip lands exactly on the target EMU
# values instead of drifting by +/-1 EMU.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

$targetId = 20
$shape = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $candidate = $s.Shapes.Item($i)
    if ($candidate.Id -eq $targetId) {
        $shape = $candidate
        break
    }
}

if ($shape -eq $null) {
    throw "Could not find shape with Id $targetId on slide 2"
}

$shape.Left = 591.7890013779527
$shape.Top = 400.4109041417323
